$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04678466666666666
$ws.Range("H2").Value = 0.140354
$ws.Range("I2").Value = 0.006739448717762189
$ws.Range("J2").Value = 0.006739448717762188
$ws.Range("M2").Value = 139.728498
$ws.Range("N2").Value = 419.185494
$ws.Range("O2").Value = 0.9065295391216045
$ws.Range("P2").Value = 0.9065295391216045
$ws.Range("Q2").Value = 6.537151202763999
$ws.Range("R2").Value = 58.83436082487599
$ws.Range("S2").Value = 0.006109509340046645
$ws.Range("T2").Value = 0.006109509340046644

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04678466666666666
$ws.Range("H3").Value = 0.140354
$ws.Range("I3").Value = 0.006739448717762189
$ws.Range("J3").Value = 0.006739448717762188
$ws.Range("O3").Value = 0.005362677585431591
$ws.Range("P3").Value = 0.005362677585431591
$ws.Range("Q3").Value = 0.03867125417844445
$ws.Range("R3").Value = 0.348041287606
$ws.Range("S3").Value = 0.00003614149057690896
$ws.Range("T3").Value = 0.00003614149057690896

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04678466666666666
$ws.Range("H4").Value = 0.140354
$ws.Range("I4").Value = 0.006739448717762189
$ws.Range("J4").Value = 0.006739448717762188
$ws.Range("O4").Value = 0.0881077832929639
$ws.Range("P4").Value = 0.0881077832929639
$ws.Range("Q4").Value = 0.6353614269255554
$ws.Range("R4").Value = 5.718252842329999
$ws.Range("S4").Value = 0.0005937978871386343
$ws.Range("T4").Value = 0.0005937978871386342

# Row 5
$ws.Range("I5").Value = 0.9176013393810419
$ws.Range("J5").Value = 0.9176013393810418
$ws.Range("M5").Value = 139.728498
$ws.Range("N5").Value = 419.185494
$ws.Range("O5").Value = 0.9065295391216045
$ws.Range("P5").Value = 0.9065295391216045
$ws.Range("Q5").Value = 890.057770390516
$ws.Range("R5").Value = 8010.519933514644
$ws.Range("S5").Value = 0.831832719286463
$ws.Range("T5").Value = 0.8318327192864629

# Row 6
$ws.Range("I6").Value = 0.9176013393810419
$ws.Range("J6").Value = 0.9176013393810418
$ws.Range("O6").Value = 0.005362677585431591
$ws.Range("P6").Value = 0.005362677585431591
$ws.Range("S6").Value = 0.00492080013506072
$ws.Range("T6").Value = 0.004920800135060719

# Row 7
$ws.Range("I7").Value = 0.9176013393810419
$ws.Range("J7").Value = 0.9176013393810418
$ws.Range("O7").Value = 0.0881077832929639
$ws.Range("P7").Value = 0.0881077832929639
$ws.Range("S7").Value = 0.08084781995951826
$ws.Range("T7").Value = 0.08084781995951826

# Row 8
$ws.Range("G8").Value = 0.5252196666666668
$ws.Range("I8").Value = 0.07565921190119594
$ws.Range("J8").Value = 0.07565921190119593
$ws.Range("M8").Value = 139.728498
$ws.Range("N8").Value = 419.185494
$ws.Range("O8").Value = 0.9065295391216045
$ws.Range("P8").Value = 0.9065295391216045
$ws.Range("Q8").Value = 73.38815514339402
$ws.Range("R8").Value = 660.4933962905461
$ws.Range("S8").Value = 0.06858731049509498
$ws.Range("T8").Value = 0.06858731049509496

# Row 9
$ws.Range("G9").Value = 0.5252196666666668
$ws.Range("I9").Value = 0.07565921190119594
$ws.Range("J9").Value = 0.07565921190119593
$ws.Range("O9").Value = 0.005362677585431591
$ws.Range("P9").Value = 0.005362677585431591
$ws.Range("R9").Value = 3.907223073001001
$ws.Range("S9").Value = 0.0004057359597939625
$ws.Range("T9").Value = 0.0004057359597939625

# Row 10
$ws.Range("G10").Value = 0.5252196666666668
$ws.Range("I10").Value = 0.07565921190119594
$ws.Range("J10").Value = 0.07565921190119593
$ws.Range("O10").Value = 0.0881077832929639
$ws.Range("P10").Value = 0.0881077832929639
$ws.Range("Q10").Value = 7.13277106878389
$ws.Range("R10").Value = 64.19493961905501
$ws.Range("S10").Value = 0.006666165446307007
$ws.Range("T10").Value = 0.006666165446307006
